$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item(2)

# Insert two blank columns before column A (pushes old A..F to C..H)
$ws2.Columns("A:B").Insert()

# Insert a new row before row 16 (old T12 "Form - submission" row)
$ws2.Rows("16:16").Insert()
